# TradingModel - 2021/11/12 data update
# Update the "Total Open Position" data table (Stock_Id / PositionSize) with
# the latest figures, shifting existing rows down and appending new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (Index, Stock_Id, PositionSize) for rows 3..10 (A3:C10).
# Row 2 (index 0 / 1711 / 220) stays unchanged.
$data = @(
    @(1, 2436, 50),
    @(2, 3033, 180),
    @(3, 3035, 32),
    @(4, 3141, 27),
    @(5, 3189, 27),
    @(7, 3588, 35),
    @(9, 6104, 36),
    @(11, 6411, 26)
)

$row = 3
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row = $row + 1
}

# Rows 9 and 10 are brand new - copy the formatting (bold/border/centered
# style) used by the rest of column A (A2:A8) onto the new A9/A10 cells.
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A8").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
